$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - moonshotai/kimi-k2-instruct
$ws.Range("I9").Value = 2
$ws.Range("J9").Value = 0.002
$ws.Range("K9").Value = 307
$ws.Range("L9").Value = 0.001023333333333333

# Row 11 - openai/gpt-oss-120b
$ws.Range("I11").Value = 2
$ws.Range("J11").Value = 0.002
$ws.Range("K11").Value = 286
$ws.Range("L11").Value = 0.00143

# Row 12 - openai/gpt-oss-20b
$ws.Range("I12").Value = 2
$ws.Range("J12").Value = 0.002
$ws.Range("K12").Value = 519
$ws.Range("L12").Value = 0.002595
